$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure target cells are treated as plain text so values like
# "97.30" or "0.09700" keep their exact formatting instead of
# being parsed/normalized as numbers.
$targetCells = @("D2","E2","D3","E3","E4","D5","E5","E6","D7","E7","D8","E8","D9","E9","D10","E10","D11","E11","D12","E12","D13","E13","D14","E14","D15","E15","D16","E16","D17","E17","D18","E18","E19","D20","E20","D21","E21","D22","E22","D23","E23","D24","E24","D25","E25","D26","E26","D27","E27","D28","E28","E29","D30","E30","D31","E31","D32","E32","D33","E33","D34","E34","D35","E35","D36","E36","D37","E37","D38","E38","D39","E39","E40","D41","E41","D42","E42","D43","E43","D44","E44","E45","D46","E46","D47","E47","D48","E48","D49","E49","D50","E50","E51")
foreach ($cell in $targetCells) {
    $ws.Range($cell).NumberFormat = "@"
}

$ws.Range("D2").Value = "30.366.51"
$ws.Range("E2").Value = "  +0.14%  "
$ws.Range("D3").Value = "1.870.75"
$ws.Range("E3").Value = "  -0.71%  "
$ws.Range("E4").Value = "  +0.28%  "
$ws.Range("D5").Value = "235.78"
$ws.Range("E5").Value = "  -0.89%  "
$ws.Range("E6").Value = "  +0.21%  "
$ws.Range("D7").Value = "0.4667"
$ws.Range("E7").Value = "  -0.13%  "
$ws.Range("D8").Value = "0.2841"
$ws.Range("E8").Value = "  +0.97%  "
$ws.Range("D9").Value = "0.06555"
$ws.Range("E9").Value = "  -0.02%  "
$ws.Range("D10").Value = "21.09"
$ws.Range("E10").Value = "  +7.68%  "
$ws.Range("D11").Value = "0.07938"
$ws.Range("E11").Value = "  +2.67%  "
$ws.Range("D12").Value = "97.30"
$ws.Range("E12").Value = "  -1.23%  "
$ws.Range("D13").Value = "1.868.75"
$ws.Range("E13").Value = "  -0.87%  "
$ws.Range("D14").Value = "5.155"
$ws.Range("E14").Value = "  +0.60%  "
$ws.Range("D15").Value = "0.6746"
$ws.Range("E15").Value = "  +0.89%  "
$ws.Range("D16").Value = "281.95"
$ws.Range("E16").Value = "  -1.31%  "
$ws.Range("D17").Value = "30.372.21"
$ws.Range("E17").Value = "  +0.21%  "
$ws.Range("D18").Value = "5.544"
$ws.Range("E18").Value = "  +4.41%  "
$ws.Range("E19").Value = "  +0.17%  "
$ws.Range("D20").Value = "12.69"
$ws.Range("E20").Value = "  +0.96%  "
$ws.Range("D21").Value = "2.116.30"
$ws.Range("E21").Value = "  -0.52%  "
$ws.Range("D22").Value = "0.000007287"
$ws.Range("E22").Value = "  -0.06%  "
$ws.Range("D23").Value = "1.002"
$ws.Range("E23").Value = "  +0.37%  "
$ws.Range("D24").Value = "6.206"
$ws.Range("E24").Value = "  +0.23%  "
$ws.Range("D25").Value = "9.292"
$ws.Range("E25").Value = "  +0.29%  "
$ws.Range("D26").Value = "164.65"
$ws.Range("E26").Value = "  -1.69%  "
$ws.Range("D27").Value = "19.08"
$ws.Range("E27").Value = "  +0.13%  "
$ws.Range("D28").Value = "1.935"
$ws.Range("E28").Value = "  -2.50%  "
$ws.Range("E29").Value = "  -1.17%  "
$ws.Range("D30").Value = "0.09700"
$ws.Range("E30").Value = "  -1.50%  "
$ws.Range("D31").Value = "4.442"
$ws.Range("E31").Value = "  -0.44%  "
$ws.Range("D32").Value = "1.478"
$ws.Range("E32").Value = "  -0.97%  "
$ws.Range("D33").Value = "4.119"
$ws.Range("E33").Value = "  -1.70%  "
$ws.Range("D34").Value = "0.04708"
$ws.Range("E34").Value = "  +0.75%  "
$ws.Range("D35").Value = "1.120"
$ws.Range("E35").Value = "  +2.18%  "
$ws.Range("D36").Value = "0.7049"
$ws.Range("E36").Value = "  -0.56%  "
$ws.Range("D37").Value = "2.718"
$ws.Range("E37").Value = "  +0.56%  "
$ws.Range("D38").Value = "0.01860"
$ws.Range("E38").Value = "  -0.57%  "
$ws.Range("D39").Value = "6.339"
$ws.Range("E39").Value = "  -5.05%  "
$ws.Range("E40").Value = "  +1.13%  "
$ws.Range("D41").Value = "73.67"
$ws.Range("E41").Value = "  +1.89%  "
$ws.Range("D42").Value = "1.949"
$ws.Range("E42").Value = "  -0.63%  "
$ws.Range("D43").Value = "0.8478"
$ws.Range("E43").Value = "  -2.42%  "
$ws.Range("D44").Value = "0.4198"
$ws.Range("E44").Value = "  +0.24%  "
$ws.Range("E45").Value = "  +0.28%  "
$ws.Range("D46").Value = "103.90"
$ws.Range("E46").Value = "  -0.05%  "
$ws.Range("D47").Value = "7.218"
$ws.Range("E47").Value = "  -0.39%  "
$ws.Range("D48").Value = "9.239"
$ws.Range("E48").Value = "  -1.47%  "
$ws.Range("D49").Value = "939.61"
$ws.Range("E49").Value = "  -5.34%  "
$ws.Range("D50").Value = "34.19"
$ws.Range("E50").Value = "  +0.47%  "
$ws.Range("E51").Value = "  -2.44%  "
